{"js": "const replacements = [\n  { find: \"2024-03-21 Thursday\", replace: \"2024-03-22 Friday\" },\n  { find: \"67\u00d762=\", replace: \"92\u00d790=\" },\n  { find: \"78\u00d763=\", replace: \"12\u00d793=\" },\n  { find: \"38\u00d766=\", replace: \"63\u00d787=\" },\n  { find: \"63\u00d782=\", replace: \"61\u00d734=\" },\n  { find: \"90\u00d776=\", replace: \"62\u00d791=\" },\n  { find: \"77\u00d738=\", replace: \"37\u00d736=\" },\n  { find: \"75\u00d767=\", replace: \"14\u00d768=\" },\n  { find: \"89\u00d713=\", replace: \"15\u00d790=\" },\n  { find: \"56\u00d799=\", replace: \"59\u00d740=\" },\n  { find: \"64\u00d775=\", replace: \"83\u00d766=\" },\n  { find: \"14\u00d755=\", replace: \"52\u00d770=\" },\n  { find: \"73\u00d717=\", replace: \"53\u00d747=\" },\n  { find: \"66\u00d752=\", replace: \"60\u00d797=\" },\n  { find: \"12\u00d764=\", replace: \"14\u00d777=\" },\n  { find: \"73\u00d789=\", replace: \"29\u00d781=\" },\n  { find: \"96\u00d773=\", replace: \"74\u00d752=\" },\n  { find: \"65\u00d798=\", replace: \"31\u00d790=\" },\n  { find: \"19\u00d750=\", replace: \"55\u00d763=\" },\n  { find: \"51\u00d792=\", replace: \"90\u00d768=\" },\n  { find: \"24\u00d787=\", replace: \"59\u00d773=\" },\n  { find: \"79\u00d737=\", replace: \"18\u00d773=\" },\n  { find: \"78\u00d748=\", replace: \"94\u00d717=\" },\n  { find: \"75\u00d735=\", replace: \"63\u00d769=\" },\n  { find: \"29\u00d761=\", replace: \"94\u00d769=\" },\n  { find: \"82\u00d779=\", replace: \"47\u00d721=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Find = '2024-03-21 Thursday'; Replace = '2024-03-22 Friday'},\n    @{Find = '67\u00d762='; Replace = '92\u00d790='},\n    @{Find = '78\u00d763='; Replace = '12\u00d793='},\n    @{Find = '38\u00d766='; Replace = '63\u00d787='},\n    @{Find = '63\u00d782='; Replace = '61\u00d734='},\n    @{Find = '90\u00d776='; Replace = '62\u00d791='},\n    @{Find = '77\u00d738='; Replace = '37\u00d736='},\n    @{Find = '75\u00d767='; Replace = '14\u00d768='},\n    @{Find = '89\u00d713='; Replace = '15\u00d790='},\n    @{Find = '56\u00d799='; Replace = '59\u00d740='},\n    @{Find = '64\u00d775='; Replace = '83\u00d766='},\n    @{Find = '14\u00d755='; Replace = '52\u00d770='},\n    @{Find = '73\u00d717='; Replace = '53\u00d747='},\n    @{Find = '66\u00d752='; Replace = '60\u00d797='},\n    @{Find = '12\u00d764='; Replace = '14\u00d777='},\n    @{Find = '73\u00d789='; Replace = '29\u00d781='},\n    @{Find = '96\u00d773='; Replace = '74\u00d752='},\n    @{Find = '65\u00d798='; Replace = '31\u00d790='},\n    @{Find = '19\u00d750='; Replace = '55\u00d763='},\n    @{Find = '51\u00d792='; Replace = '90\u00d768='},\n    @{Find = '24\u00d787='; Replace = '59\u00d773='},\n    @{Find = '79\u00d737='; Replace = '18\u00d773='},\n    @{Find = '78\u00d748='; Replace = '94\u00d717='},\n    @{Find = '75\u00d735='; Replace = '63\u00d769='},\n    @{Find = '29\u00d761='; Replace = '94\u00d769='},\n    @{Find = '82\u00d779='; Replace = '47\u00d721='}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.Replace\n    $find.Execute(\n        $r.Find,        # FindText\n        $true,          # MatchCase\n        $true,          # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $r.Replace,     # ReplaceWith\n        2               # Replace (wdReplaceAll)\n    )\n}\n"}
